# Updates currentAveragePrice / LevePrice / LeveProfit figures that were
# refreshed by the scheduled price-sync runner across several Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20: Shut Up and Take My Gil
$ws.Range("H20").Value = 5001500
$ws.Range("I20").Value = 5001500
$ws.Range("K20").Value = 5001500
$ws.Range("M20").Value = -5001270

# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 1229.9
$ws.Range("I28").Value = 1119.7778
$ws.Range("K28").Value = 1119.7778
$ws.Range("M28").Value = -634.7778000000001

# Row 35: Conspicuous Conjuration
$ws.Range("H35").Value = 5001500
$ws.Range("I35").Value = 5001500
$ws.Range("K35").Value = 5001500
$ws.Range("M35").Value = -5001121

# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 30669.422
$ws.Range("I40").Value = 51250
$ws.Range("J40").Value = 28248.176
$ws.Range("K40").Value = 51250
$ws.Range("L40").Value = 28248.176
$ws.Range("M40").Value = -51075
$ws.Range("N40").Value = -28598.176

# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 4208.1665
$ws.Range("I62").Value = 3785.7144
$ws.Range("J62").Value = 4799.6
$ws.Range("K62").Value = 3785.7144
$ws.Range("L62").Value = 4799.6
$ws.Range("M62").Value = -3161.7144
$ws.Range("N62").Value = -6047.6

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 4208.1665
$ws.Range("I65").Value = 3785.7144
$ws.Range("J65").Value = 4799.6
$ws.Range("K65").Value = 18928.572
$ws.Range("L65").Value = 23998
$ws.Range("M65").Value = -15808.572
$ws.Range("N65").Value = -30238

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 4998.6665
$ws.Range("I76").Value = 4998
$ws.Range("K76").Value = 4998
$ws.Range("M76").Value = -4683

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 4998.6665
$ws.Range("I79").Value = 4998
$ws.Range("K79").Value = 4998
$ws.Range("M79").Value = -3906

# Row 113: Amaro Kart
$ws.Range("H113").Value = 3599.5715
$ws.Range("I113").Value = 3456.7144
$ws.Range("K113").Value = 3456.7144
$ws.Range("M113").Value = -202.7143999999998

# Row 125: Body over Mind
$ws.Range("H125").Value = 2278.5715
$ws.Range("I125").Value = 1690
$ws.Range("K125").Value = 15210
$ws.Range("M125").Value = -12750

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 5587.7856
$ws.Range("I132").Value = 2024.1
$ws.Range("K132").Value = 6072.299999999999
$ws.Range("M132").Value = -3542.299999999999

# Row 135: For Tired Minds
$ws.Range("H135").Value = 4263.316
$ws.Range("I135").Value = 1320.5385
$ws.Range("K135").Value = 11884.8465
$ws.Range("M135").Value = -9349.846500000001

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 3124.5557
$ws.Range("I137").Value = 2817.3572
$ws.Range("K137").Value = 8452.071599999999
$ws.Range("M137").Value = -5902.071599999999

# Row 138: All-night Crafting
$ws.Range("H138").Value = 6151.8477
$ws.Range("J138").Value = 6974.1025
$ws.Range("L138").Value = 20922.3075
$ws.Range("N138").Value = -31202.3075

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 46363.273
$ws.Range("I2").Value = 50929.6
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 50929.6
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = -50816.6
$ws.Range("N2").Value = -926

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 20349.727
$ws.Range("I32").Value = 19814.342
$ws.Range("K32").Value = 19814.342
$ws.Range("M32").Value = -19527.342

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2191.8333
$ws.Range("I45").Value = 1538.3334
$ws.Range("K45").Value = 1538.3334
$ws.Range("M45").Value = -1161.3334

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 22059438
$ws.Range("I74").Value = 25862582
$ws.Range("J74").Value = 1199.6
$ws.Range("K74").Value = 25862582
$ws.Range("L74").Value = 1199.6
$ws.Range("M74").Value = -25861708
$ws.Range("N74").Value = -2947.6

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 22059438
$ws.Range("I77").Value = 25862582
$ws.Range("J77").Value = 1199.6
$ws.Range("K77").Value = 129312910
$ws.Range("L77").Value = 5998
$ws.Range("M77").Value = -129308542
$ws.Range("N77").Value = -14734

# Row 101: Art Imitates Life
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 116: No Scope
$ws.Range("H116").Value = 46363.273
$ws.Range("I116").Value = 50929.6
$ws.Range("J116").Value = 700
$ws.Range("K116").Value = 50929.6
$ws.Range("L116").Value = 700
$ws.Range("M116").Value = -48635.6
$ws.Range("N116").Value = -5288

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 46363.273
$ws.Range("I3").Value = 50929.6
$ws.Range("J3").Value = 700
$ws.Range("K3").Value = 50929.6
$ws.Range("L3").Value = 700
$ws.Range("M3").Value = -50815.6
$ws.Range("N3").Value = -928

$ws = $wb.Worksheets.Item("CRP")
# Row 88: Hold on Adamantite
$ws.Range("H88").Value = 79999
$ws.Range("J88").Value = 79999
$ws.Range("L88").Value = 79999
$ws.Range("N88").Value = -80811

# Row 91: Spears for Stone Vigilantes (L)
$ws.Range("H91").Value = 79999
$ws.Range("J91").Value = 79999
$ws.Range("L91").Value = 79999
$ws.Range("N91").Value = -82807

# Row 133: Yimepi's Country Charms
$ws.Range("H133").Value = 59655
$ws.Range("J133").Value = 59655
$ws.Range("L133").Value = 59655
$ws.Range("N133").Value = -64715

# Row 139: Weaving a Path
$ws.Range("H139").Value = 58422.9
$ws.Range("J139").Value = 94846
$ws.Range("L139").Value = 94846
$ws.Range("N139").Value = -105126

# Row 141: No Greater Treasure
$ws.Range("H141").Value = 110520.1
$ws.Range("I141").Value = 90771.336
$ws.Range("J141").Value = 111780.66
$ws.Range("K141").Value = 90771.336
$ws.Range("L141").Value = 111780.66
$ws.Range("M141").Value = -85591.336
$ws.Range("N141").Value = -122140.66

$ws = $wb.Worksheets.Item("GSM")
# Row 124: The Sage's Successor
$ws.Range("H124").Value = 361762.72
$ws.Range("J124").Value = 361762.72
$ws.Range("L124").Value = 361762.72
$ws.Range("N124").Value = -371582.72

# Row 132: On Board for Lar
$ws.Range("H132").Value = 3064.2942
$ws.Range("I132").Value = 2842.7354
$ws.Range("K132").Value = 8528.206200000001
$ws.Range("M132").Value = -5998.206200000001

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 3540.7
$ws.Range("J16").Value = 9333.666999999999
$ws.Range("L16").Value = 9333.666999999999
$ws.Range("N16").Value = -9673.666999999999

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 343.4737
$ws.Range("J55").Value = 738
$ws.Range("L55").Value = 738
$ws.Range("N55").Value = -1084

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4535.2104
$ws.Range("I136").Value = 3262
$ws.Range("K136").Value = 9786
$ws.Range("M136").Value = -7236

# Row 137: Lending Artisans a Hand
$ws.Range("H137").Value = 123456
$ws.Range("J137").Value = 123456
$ws.Range("L137").Value = 123456
$ws.Range("N137").Value = -133656

# Row 138: Freezing Toes
$ws.Range("H138").Value = 100429
$ws.Range("J138").Value = 100429
$ws.Range("L138").Value = 100429
$ws.Range("N138").Value = -110709

# Row 139: Giving Gatherers Their Gear
$ws.Range("H139").Value = 83499
$ws.Range("J139").Value = 83499
$ws.Range("L139").Value = 83499
$ws.Range("N139").Value = -93779

# Row 141: Just Generally Freezing
$ws.Range("H141").Value = 97999
$ws.Range("I141").Value = 97998
$ws.Range("J141").Value = 98000
$ws.Range("K141").Value = 97998
$ws.Range("L141").Value = 98000
$ws.Range("M141").Value = -92818
$ws.Range("N141").Value = -108360

$ws = $wb.Worksheets.Item("WVR")
# Row 26: New Shoes, New Me
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

# Row 92: Modest Beginnings
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 107: Flax Wax
$ws.Range("H107").Value = 1182.8334
$ws.Range("I107").Value = 1085.2858
$ws.Range("J107").Value = 1319.4
$ws.Range("K107").Value = 3255.8574
$ws.Range("L107").Value = 3958.2
$ws.Range("M107").Value = -1335.8574
$ws.Range("N107").Value = -7798.200000000001

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 7486.7
$ws.Range("I136").Value = 5912.8335
$ws.Range("J136").Value = 9847.5
$ws.Range("K136").Value = 17738.5005
$ws.Range("L136").Value = 29542.5
$ws.Range("M136").Value = -15188.5005
$ws.Range("N136").Value = -34642.5

# Row 138: Halfgloves, Full Effort
$ws.Range("H138").Value = 133874.75
$ws.Range("J138").Value = 148600
$ws.Range("L138").Value = 148600
$ws.Range("N138").Value = -158880

# Row 140: Glamorous Gloves
$ws.Range("H140").Value = 79668.75
$ws.Range("J140").Value = 79668.75
$ws.Range("L140").Value = 79668.75
$ws.Range("N140").Value = -90028.75
